$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: Remi Bossard
$ws.Range("A23").Value = "Remi Bossard"
$ws.Range("B21").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "BA(Mod) Molecular Medicine"
$ws.Range("C23").Value = 2025
$ws.Range("D23").Value = "Trinity College Dublin"

# Row 24: Linda Romano
$ws.Range("A24").Value = "Linda Romano"
$ws.Range("B21").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").Value = "BA(Mod) Molecular Medicine"
$ws.Range("C24").Value = 2025
$ws.Range("D24").Value = "Trinity College Dublin"

$ws.Range("D28").Select() | Out-Null
